# Update the classification-report text (column A) with the new
# precision/recall/f1/accuracy figures from the run without StandardScaler.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "           0       0.98      0.93      0.95      7115"
$ws.Range("A4").Value = "           1       0.93      0.98      0.96      7114"
$ws.Range("A6").Value = "    accuracy                           0.95     14229"
$ws.Range("A7").Value = "   macro avg       0.96      0.95      0.95     14229"
$ws.Range("A8").Value = "weighted avg       0.96      0.95      0.95     14229"

# Row 13 holds the per-epoch training-accuracy list as text cells
# (A13 starts with "[", J13 ends with "]"). B13:I13 look like plain
# numbers, so they need a leading apostrophe to keep Excel from storing
# them as numeric values - they must stay text, same as before the edit.
$ws.Range("A13").Value = "[0.8724017143249512"
$ws.Range("B13").Value = "'" + " 0.9084040522575378"
$ws.Range("C13").Value = "'" + " 0.9285224676132202"
$ws.Range("D13").Value = "'" + " 0.9432643055915833"
$ws.Range("E13").Value = "'" + " 0.9518036246299744"
$ws.Range("F13").Value = "'" + " 0.9595522880554199"
$ws.Range("G13").Value = "'" + " 0.9646302461624146"
$ws.Range("H13").Value = "'" + " 0.9695851802825928"
$ws.Range("I13").Value = "'" + " 0.9722734689712524"
$ws.Range("J13").Value = " 0.9751023650169373]"

# Drop the auto-applied "quote prefix" style so the cells keep their
# original (unstyled) formatting - only the text content changed.
$ws.Range("B13:I13").Style = "Normal"
